# Octobre 2023 update: append daily rows 165-195 (dates 2023-10-01..2023-10-31)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A on the new rows uses the same custom date/time format as the rest of column A
$ws.Range("A165:A195").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 165
$ws.Range("A165").Value = 45200
$ws.Range("D165").Value = 65850
$ws.Range("F165").Value = 123500
$ws.Range("G165").Value = 2027450
$ws.Range("I165").Value = 14000
$ws.Range("J165").Value = 1000
$ws.Range("R165").Value = 5000

# Row 166
$ws.Range("A166").Value = 45201
$ws.Range("C166").Value = 77500
$ws.Range("D166").Value = 215525
$ws.Range("J166").Value = 1000
$ws.Range("R166").Value = 9000
$ws.Range("S166").Value = 4400
$ws.Range("W166").Value = 40000

# Row 167
$ws.Range("A167").Value = 45202
$ws.Range("D167").Value = 19000
$ws.Range("F167").Value = 209000
$ws.Range("I167").Value = 10500
$ws.Range("J167").Value = 1000
$ws.Range("V167").Value = 247500

# Row 168
$ws.Range("A168").Value = 45203
$ws.Range("D168").Value = 79350
$ws.Range("I168").Value = 7000
$ws.Range("J168").Value = 1000
$ws.Range("K168").Value = 50000
$ws.Range("M168").Value = 250000
$ws.Range("Q168").Value = 1500
$ws.Range("R168").Value = 10600

# Row 169
$ws.Range("A169").Value = 45204
$ws.Range("C169").Value = 77500
$ws.Range("D169").Value = 121550
$ws.Range("F169").Value = 886000
$ws.Range("I169").Value = 5000
$ws.Range("J169").Value = 1000
$ws.Range("K169").Value = 55000
$ws.Range("V169").Value = 1000

# Row 170
$ws.Range("A170").Value = 45205
$ws.Range("D170").Value = 99750
$ws.Range("F170").Value = 187000
$ws.Range("I170").Value = 10500
$ws.Range("J170").Value = 1000
$ws.Range("R170").Value = 42000
$ws.Range("U170").Value = 3700

# Row 171
$ws.Range("A171").Value = 45206
$ws.Range("D171").Value = 112100
$ws.Range("F171").Value = 118400
$ws.Range("I171").Value = 8750
$ws.Range("J171").Value = 1000
$ws.Range("R171").Value = 30000
$ws.Range("S171").Value = 600
$ws.Range("V171").Value = 12500

# Row 172
$ws.Range("A172").Value = 45207
$ws.Range("C172").Value = 77500
$ws.Range("D172").Value = 54450
$ws.Range("F172").Value = 295500
$ws.Range("G172").Value = 1856550
$ws.Range("I172").Value = 17500
$ws.Range("J172").Value = 1000
$ws.Range("S172").Value = 160000
$ws.Range("T172").Value = 500
$ws.Range("V172").Value = 435000

# Row 173
$ws.Range("A173").Value = 45208
$ws.Range("C173").Value = 15000
$ws.Range("D173").Value = 106950
$ws.Range("E173").Value = 37250
$ws.Range("F173").Value = 30000
$ws.Range("J173").Value = 1000
$ws.Range("R173").Value = 10600
$ws.Range("U173").Value = 5000

# Row 174
$ws.Range("A174").Value = 45209
$ws.Range("D174").Value = 17500
$ws.Range("F174").Value = 41000
$ws.Range("I174").Value = 10500
$ws.Range("J174").Value = 1000
$ws.Range("R174").Value = 1600

# Row 175
$ws.Range("A175").Value = 45210
$ws.Range("C175").Value = 77500
$ws.Range("D175").Value = 58125
$ws.Range("F175").Value = 620500
$ws.Range("I175").Value = 10500
$ws.Range("J175").Value = 1000
$ws.Range("K175").Value = 90000
$ws.Range("R175").Value = 10600
$ws.Range("U175").Value = 90000
$ws.Range("V175").Value = 35000

# Row 176
$ws.Range("A176").Value = 45211
$ws.Range("D176").Value = 163700
$ws.Range("F176").Value = 126400
$ws.Range("I176").Value = 10500
$ws.Range("J176").Value = 1000
$ws.Range("K176").Value = 35000
$ws.Range("Q176").Value = 1400
$ws.Range("R176").Value = 10600

# Row 177
$ws.Range("A177").Value = 45212
$ws.Range("D177").Value = 149650
$ws.Range("F177").Value = 393200
$ws.Range("J177").Value = 1000
$ws.Range("T177").Value = 2000
$ws.Range("U177").Value = 157800

# Row 178
$ws.Range("A178").Value = 45213
$ws.Range("C178").Value = 77500
$ws.Range("D178").Value = 82600
$ws.Range("F178").Value = 82000
$ws.Range("J178").Value = 1000
$ws.Range("M178").Value = 250000
$ws.Range("R178").Value = 9000
$ws.Range("U178").Value = 80700
$ws.Range("V178").Value = 170000

# Row 179
$ws.Range("A179").Value = 45214
$ws.Range("C179").Value = 77500
$ws.Range("D179").Value = 155100
$ws.Range("F179").Value = 134000
$ws.Range("G179").Value = 1440000
$ws.Range("I179").Value = 31500
$ws.Range("J179").Value = 1000
$ws.Range("R179").Value = 9000

# Row 180
$ws.Range("A180").Value = 45215
$ws.Range("D180").Value = 141650
$ws.Range("J180").Value = 1000
$ws.Range("O180").Value = 55000
$ws.Range("R180").Value = 39500
$ws.Range("T180").Value = 5000
$ws.Range("W180").Value = 15000

# Row 181
$ws.Range("A181").Value = 45216
$ws.Range("D181").Value = 19000
$ws.Range("F181").Value = 336000
$ws.Range("J181").Value = 1000
$ws.Range("Q181").Value = 2000
$ws.Range("R181").Value = 32000

# Row 182
$ws.Range("A182").Value = 45217
$ws.Range("D182").Value = 100050
$ws.Range("F182").Value = 65400
$ws.Range("I182").Value = 10000
$ws.Range("J182").Value = 1000
$ws.Range("K182").Value = 50000
$ws.Range("Q182").Value = 1000
$ws.Range("R182").Value = 4100
$ws.Range("T182").Value = 1000
$ws.Range("V182").Value = 18300

# Row 183
$ws.Range("A183").Value = 45218
$ws.Range("C183").Value = 77500
$ws.Range("D183").Value = 170400
$ws.Range("F183").Value = 73500
$ws.Range("J183").Value = 1000
$ws.Range("K183").Value = 35000
$ws.Range("R183").Value = 9000
$ws.Range("V183").Value = 43500

# Row 184
$ws.Range("A184").Value = 45219
$ws.Range("D184").Value = 68350
$ws.Range("F184").Value = 413000
$ws.Range("I184").Value = 10500
$ws.Range("J184").Value = 1000
$ws.Range("R184").Value = 10600
$ws.Range("V184").Value = 13000

# Row 185
$ws.Range("A185").Value = 45220
$ws.Range("C185").Value = 77500
$ws.Range("D185").Value = 145650
$ws.Range("I185").Value = 10500
$ws.Range("J185").Value = 1000
$ws.Range("M185").Value = 75000
$ws.Range("V185").Value = 62000

# Row 186
$ws.Range("A186").Value = 45221
$ws.Range("D186").Value = 115250
$ws.Range("F186").Value = 57000
$ws.Range("G186").Value = 1696900
$ws.Range("I186").Value = 17500
$ws.Range("J186").Value = 1000
$ws.Range("R186").Value = 2800
$ws.Range("V186").Value = 21000

# Row 187
$ws.Range("A187").Value = 45222
$ws.Range("D187").Value = 203050
$ws.Range("F187").Value = 642500
$ws.Range("J187").Value = 1000
$ws.Range("R187").Value = 64000
$ws.Range("S187").Value = 9000
$ws.Range("U187").Value = 20000

# Row 188
$ws.Range("A188").Value = 45223
$ws.Range("D188").Value = 25500
$ws.Range("F188").Value = 89000
$ws.Range("I188").Value = 10500
$ws.Range("J188").Value = 1000

# Row 189
$ws.Range("A189").Value = 45224
$ws.Range("C189").Value = 77500
$ws.Range("D189").Value = 120325
$ws.Range("F189").Value = 51500
$ws.Range("I189").Value = 10500
$ws.Range("J189").Value = 1000
$ws.Range("K189").Value = 60000
$ws.Range("R189").Value = 9000
$ws.Range("S189").Value = 500

# Row 190
$ws.Range("A190").Value = 45225
$ws.Range("D190").Value = 121800
$ws.Range("F190").Value = 858900
$ws.Range("I190").Value = 10500
$ws.Range("J190").Value = 2000
$ws.Range("K190").Value = 35000
$ws.Range("R190").Value = 10600
$ws.Range("T190").Value = 11000
$ws.Range("U190").Value = 37000

# Row 191
$ws.Range("A191").Value = 45226
$ws.Range("D191").Value = 64475
$ws.Range("E191").Value = 37250
$ws.Range("F191").Value = 38000
$ws.Range("I191").Value = 10500
$ws.Range("J191").Value = 1000
$ws.Range("M191").Value = 300000
$ws.Range("R191").Value = 14000
$ws.Range("T191").Value = 1000
$ws.Range("V191").Value = 10000

# Row 192
$ws.Range("A192").Value = 45227
$ws.Range("C192").Value = 77500
$ws.Range("D192").Value = 180700
$ws.Range("F192").Value = 31100
$ws.Range("I192").Value = 7000
$ws.Range("J192").Value = 1000
$ws.Range("R192").Value = 30000
$ws.Range("T192").Value = 1000

# Row 193
$ws.Range("A193").Value = 45228
$ws.Range("D193").Value = 102700
$ws.Range("F193").Value = 57000
$ws.Range("G193").Value = 919250
$ws.Range("I193").Value = 17500
$ws.Range("J193").Value = 1000
$ws.Range("R193").Value = 9000
$ws.Range("T193").Value = 5000
$ws.Range("W193").Value = 190000

# Row 194
$ws.Range("A194").Value = 45229
$ws.Range("D194").Value = 62800
$ws.Range("F194").Value = 47500
$ws.Range("J194").Value = 1000
$ws.Range("W194").Value = 289500

# Row 195
$ws.Range("A195").Value = 45230
$ws.Range("C195").Value = 77.5
$ws.Range("D195").Value = 25.5
$ws.Range("F195").Value = 7
$ws.Range("J195").Value = 1000
$ws.Range("P195").Value = 200000
$ws.Range("U195").Value = 42000
$ws.Range("V195").Value = 10000
